# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Swap "Peru" / "India" shared strings: the ranking reshuffled so that
# the row previously holding Peru's stats (row 16) now shows India's new
# numbers, while the old Peru numbers slide down to row 17 (still labeled
# Peru). The country label in A16 stays "India" going forward (formerly
# "Peru"), and A17 stays "Peru" (formerly "India") because the underlying
# shared-string table entries were swapped; net visible effect is that
# row 16 = India (new data) and row 17 = Peru (row16's old data).
$ws.Range("A16").Value = "India"
$ws.Range("A17").Value = "Peru"

# Row 13 (rank 17) updated values
$ws.Range("B13").Value = 109286
$ws.Range("C13").Value = 1683
$ws.Range("D13").Value = 87422
$ws.Range("E13").Value = 15179
$ws.Range("F13").Value = 2703
$ws.Range("G13").Value = 45
$ws.Range("H13").Value = 6685

# Row 16 (India) updated values
$ws.Range("B16").Value = 67700
$ws.Range("C16").Value = 539
$ws.Range("D16").Value = 21130
$ws.Range("E16").Value = 44355
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 2215

# Row 17 (Peru) updated values
$ws.Range("B17").Value = 67307
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 21349
$ws.Range("E17").Value = 44069
$ws.Range("F17").Value = 774
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 1889

# Row 39 (rank 43) updated values
$ws.Range("B39").Value = 15588
$ws.Range("C39").Value = 226
$ws.Range("D39").Value = 7245
$ws.Range("E39").Value = 7371
$ws.Range("F39").Value = 255

# Row 86 (rank 90) updated values
$ws.Range("B86").Value = 1664
$ws.Range("C86").Value = 22
$ws.Range("D86").Value = 1200
$ws.Range("E86").Value = 373

# Row 95 updated values
$ws.Range("D95").Value = 985
$ws.Range("E95").Value = 59
